$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "cryptos" price table (columns B-E, rows 2-51) with the
# latest scrape: updated Price / Volume(1h) strings for every coin, plus
# a handful of rows where two adjacent coins swapped rank order (their
# Coin name + Link + Price + Volume all move down/up one row).
#
# Price values that look like plain numbers (e.g. "0.9993", "308.59")
# are written with a leading apostrophe, the same trick Excel's UI uses
# to force text entry, so they stay text (matching how this sheet stores
# every other price, e.g. "26.680.01") instead of being auto-converted
# to a numeric value.
$ws.Range("D2").Value = "26.680.01"
$ws.Range("E2").Value = "  -1.76%  "
$ws.Range("D3").Value = "1.792.14"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "'308.59"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").Value = "'0.9988"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").Value = "'0.4460"
$ws.Range("E7").Value = "  +5.32%  "
$ws.Range("D8").Value = "'0.3657"
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("D9").Value = "'0.07301"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").Value = "'0.8554"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "2.032.04"
$ws.Range("E11").Value = "  +11.37%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'20.56"
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("D13").Value = "'6.592"
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'92.01"
$ws.Range("E14").Value = "  +2.79%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "'0.07073"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "'5.260"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "'0.9999"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "'0.000008644"
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("D19").Value = "'0.9993"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "'14.77"
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("D21").Value = "26.737.69"
$ws.Range("E21").Value = "  -1.86%  "
$ws.Range("D22").Value = "'5.128"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "'10.76"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("D24").Value = "'1.978"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'151.66"
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'18.38"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.165"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").Value = "'5.182"
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("D29").Value = "'117.12"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").Value = "'0.08792"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").Value = "'0.7397"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").Value = "'1.154"
$ws.Range("E32").Value = "  -2.94%  "
$ws.Range("D33").Value = "'2.910"
$ws.Range("E33").Value = "  -3.72%  "
$ws.Range("D34").Value = "'4.439"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "'0.9978"
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("D36").Value = "'1.083"
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("D37").Value = "'0.01952"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").Value = "'0.05158"
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("D39").Value = "'0.5299"
$ws.Range("E39").Value = "  +5.32%  "
$ws.Range("D40").Value = "'2.834"
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("D41").Value = "'7.009"
$ws.Range("E41").Value = "  -2.92%  "
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("D43").Value = "'0.5079"
$ws.Range("E43").Value = "  +7.39%  "
$ws.Range("D44").Value = "'8.379"
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("D45").Value = "'10.47"
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("D46").Value = "'1.943"
$ws.Range("E46").Value = "  +4.16%  "
$ws.Range("D47").Value = "'105.02"
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("D48").Value = "'0.9975"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("D49").Value = "'1.661"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'0.06292"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("D51").Value = "'0.9123"
$ws.Range("E51").Value = "  +0.36%  "
